$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Activate()

# Row 17 result was NG, fix it to OK
$ws.Range("E17").Value = "OK"

# New row 18: new test case, still NG
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "криво показывается каунтер у повтора: I вместо i/j"
$ws.Range("E18").Value = "NG"

# New row 19: another new test case, still NG
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "выполнили команду, увеличили каунтер -- выполняется, но кривой каунтер"
$ws.Range("E19").Value = "NG"

# Rows 20-30: just numbering continues in column A
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23
$ws.Range("A25").Value = 24
$ws.Range("A26").Value = 25
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29

# Update selection to match target state
$ws.Range("G17").Select() | Out-Null
